$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 8884
$ws.Range("E3").Value = 14450
$ws.Range("E4").Value = 6487
$ws.Range("E5").Value = 3949
$ws.Range("E6").Value = 8646
$ws.Range("E7").Value = 17249
$ws.Range("E8").Value = 17408
$ws.Range("E9").Value = 13802
$ws.Range("E10").Value = 19865
$ws.Range("E11").Value = 14936
$ws.Range("E12").Value = 3919
$ws.Range("E13").Value = 10156
